# Generate Report for Handback
# Update the handoff/handback timestamp cells to reflect the newly generated
# report times.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: "Latest HO Xliff Generate Date" for the first file row ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-08-30 03:08:42"

# --- zh-cn sheet: handoff / handback datetimes for the first file row ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-08-30 03:08:37"
$wsZhCn.Range("K2").Value = "2016-08-30 03:08:54"

# --- de-de sheet: handoff / handback datetimes for the first file row ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-08-30 03:08:42"
$wsDeDe.Range("K2").Value = "2016-08-30 03:09:03"
